$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")
Write-Host "Sheet found: $($ws.Name)"
